$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.97584633333333
$ws.Range("H2").Value = 59.927539
$ws.Range("I2").Value = 0.5797382655268795
$ws.Range("J2").Value = 0.6587494728716036
$ws.Range("M2").Value = 19.97584633333333
$ws.Range("N2").Value = 59.927539
$ws.Range("O2").Value = 0.5797382655268795
$ws.Range("P2").Value = 0.6587494728716036
$ws.Range("Q2").Value = 399.0344367329467
$ws.Range("R2").Value = 3591.309930596521
$ws.Range("S2").Value = 0.3360964565161146
$ws.Range("T2").Value = 0.4339508680086155
$ws.Range("G3").Value = 19.97584633333333
$ws.Range("H3").Value = 59.927539
$ws.Range("I3").Value = 0.5797382655268795
$ws.Range("J3").Value = 0.6587494728716036
$ws.Range("O3").Value = 0.05975306730825326
$ws.Range("P3").Value = 0.06789667671151554
$ws.Range("Q3").Value = 41.12809689169845
$ws.Range("R3").Value = 370.152872025286
$ws.Range("S3").Value = 0.03464113960119763
$ws.Range("T3").Value = 0.04472689999344455
$ws.Range("G4").Value = 19.97584633333333
$ws.Range("H4").Value = 59.927539
$ws.Range("I4").Value = 0.5797382655268795
$ws.Range("J4").Value = 0.6587494728716036
$ws.Range("M4").Value = 0.01061633333333333
$ws.Range("N4").Value = 0.031849
$ws.Range("O4").Value = 0.0003081068291285177
$ws.Range("P4").Value = 0.0003500980068860779
$ws.Range("Q4").Value = 0.2120702432901111
$ws.Range("R4").Value = 1.908632189611
$ws.Range("S4").Value = 0.0001786213187159535
$ws.Range("T4").Value = 0.0002306268774896029
$ws.Range("G5").Value = 19.97584633333333
$ws.Range("H5").Value = 59.927539
$ws.Range("I5").Value = 0.5797382655268795
$ws.Range("J5").Value = 0.6587494728716036
$ws.Range("M5").Value = 12.3983215
$ws.Range("N5").Value = 24.796643
$ws.Range("O5").Value = 0.3598236230852706
$ws.Range("P5").Value = 0.2725754432404665
$ws.Range("Q5").Value = 247.6669650752628
$ws.Range("R5").Value = 1486.001790451577
$ws.Range("S5").Value = 0.2086035231430524
$ws.Range("T5").Value = 0.179558929552401
$ws.Range("G6").Value = 19.97584633333333
$ws.Range("H6").Value = 59.927539
$ws.Range("I6").Value = 0.5797382655268795
$ws.Range("J6").Value = 0.6587494728716036
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.012988
$ws.Range("N6").Value = 0.038964
$ws.Range("O6").Value = 0.0003769372504682585
$ws.Range("P6").Value = 0.0004283091695283727
$ws.Range("Q6").Value = 0.2594462921773333
$ws.Range("R6").Value = 2.335016629596
$ws.Range("S6").Value = 0.0002185249477989391
$ws.Range("T6").Value = 0.0002821484396528898
$ws.Range("I7").Value = 0.05975306730825326
$ws.Range("J7").Value = 0.06789667671151554
$ws.Range("M7").Value = 19.97584633333333
$ws.Range("N7").Value = 59.927539
$ws.Range("O7").Value = 0.5797382655268795
$ws.Range("P7").Value = 0.6587494728716036
$ws.Range("Q7").Value = 41.12809689169845
$ws.Range("R7").Value = 370.152872025286
$ws.Range("S7").Value = 0.03464113960119763
$ws.Range("T7").Value = 0.04472689999344455
$ws.Range("I8").Value = 0.05975306730825326
$ws.Range("J8").Value = 0.06789667671151554
$ws.Range("O8").Value = 0.05975306730825326
$ws.Range("P8").Value = 0.06789667671151554
$ws.Range("S8").Value = 0.003570429052744645
$ws.Range("T8").Value = 0.004609958708468058
$ws.Range("I9").Value = 0.05975306730825326
$ws.Range("J9").Value = 0.06789667671151554
$ws.Range("M9").Value = 0.01061633333333333
$ws.Range("N9").Value = 0.031849
$ws.Range("O9").Value = 0.0003081068291285177
$ws.Range("P9").Value = 0.0003500980068860779
$ws.Range("Q9").Value = 0.02185787669177778
$ws.Range("R9").Value = 0.196720890226
$ws.Range("S9").Value = 0.00001841032809904881
$ws.Range("T9").Value = 0.00002377049119088998
$ws.Range("I10").Value = 0.05975306730825326
$ws.Range("J10").Value = 0.06789667671151554
$ws.Range("M10").Value = 12.3983215
$ws.Range("N10").Value = 24.796643
$ws.Range("O10").Value = 0.3598236230852706
$ws.Range("P10").Value = 0.2725754432404665
$ws.Range("Q10").Value = 25.52679668423033
$ws.Range("R10").Value = 153.160780105382
$ws.Range("S10").Value = 0.02150056516931373
$ws.Range("T10").Value = 0.01850696674919601
$ws.Range("I11").Value = 0.05975306730825326
$ws.Range("J11").Value = 0.06789667671151554
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.012988
$ws.Range("N11").Value = 0.038964
$ws.Range("O11").Value = 0.0003769372504682585
$ws.Range("P11").Value = 0.0004283091695283727
$ws.Range("Q11").Value = 0.02674088063733333
$ws.Range("R11").Value = 0.240667925736
$ws.Range("S11").Value = 0.00002252315689821777
$ws.Range("T11").Value = 0.00002908076921604562
$ws.Range("G12").Value = 0.01061633333333333
$ws.Range("H12").Value = 0.031849
$ws.Range("I12").Value = 0.0003081068291285177
$ws.Range("J12").Value = 0.0003500980068860779
$ws.Range("M12").Value = 19.97584633333333
$ws.Range("N12").Value = 59.927539
$ws.Range("O12").Value = 0.5797382655268795
$ws.Range("P12").Value = 0.6587494728716036
$ws.Range("Q12").Value = 0.2120702432901111
$ws.Range("R12").Value = 1.908632189611
$ws.Range("S12").Value = 0.0001786213187159535
$ws.Range("T12").Value = 0.0002306268774896029
$ws.Range("G13").Value = 0.01061633333333333
$ws.Range("H13").Value = 0.031849
$ws.Range("I13").Value = 0.0003081068291285177
$ws.Range("J13").Value = 0.0003500980068860779
$ws.Range("O13").Value = 0.05975306730825326
$ws.Range("P13").Value = 0.06789667671151554
$ws.Range("Q13").Value = 0.02185787669177778
$ws.Range("R13").Value = 0.196720890226
$ws.Range("S13").Value = 0.00001841032809904881
$ws.Range("T13").Value = 0.00002377049119088998
$ws.Range("G14").Value = 0.01061633333333333
$ws.Range("H14").Value = 0.031849
$ws.Range("I14").Value = 0.0003081068291285177
$ws.Range("J14").Value = 0.0003500980068860779
$ws.Range("M14").Value = 0.01061633333333333
$ws.Range("N14").Value = 0.031849
$ws.Range("O14").Value = 0.0003081068291285177
$ws.Range("P14").Value = 0.0003500980068860779
$ws.Range("Q14").Value = 0.0001127065334444444
$ws.Range("R14").Value = 0.001014358801
$ws.Range("S14").Value = 0.00000009492981815562963
$ws.Range("T14").Value = 0.0000001225686144256043
$ws.Range("G15").Value = 0.01061633333333333
$ws.Range("H15").Value = 0.031849
$ws.Range("I15").Value = 0.0003081068291285177
$ws.Range("J15").Value = 0.0003500980068860779
$ws.Range("M15").Value = 12.3983215
$ws.Range("N15").Value = 24.796643
$ws.Range("O15").Value = 0.3598236230852706
$ws.Range("P15").Value = 0.2725754432404665
$ws.Range("Q15").Value = 0.1316247138178333
$ws.Range("R15").Value = 0.789748282907
$ws.Range("S15").Value = 0.0001108641155543376
$ws.Range("T15").Value = 0.00009542811940457657
$ws.Range("G16").Value = 0.01061633333333333
$ws.Range("H16").Value = 0.031849
$ws.Range("I16").Value = 0.0003081068291285177
$ws.Range("J16").Value = 0.0003500980068860779
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.012988
$ws.Range("N16").Value = 0.038964
$ws.Range("O16").Value = 0.0003769372504682585
$ws.Range("P16").Value = 0.0004283091695283727
$ws.Range("Q16").Value = 0.0001378849373333333
$ws.Range("R16").Value = 0.001240964436
$ws.Range("S16").Value = 0.000000116136941022197
$ws.Range("T16").Value = 0.0000001499501865829145
$ws.Range("G17").Value = 12.3983215
$ws.Range("H17").Value = 24.796643
$ws.Range("I17").Value = 0.3598236230852706
$ws.Range("J17").Value = 0.2725754432404665
$ws.Range("M17").Value = 19.97584633333333
$ws.Range("N17").Value = 59.927539
$ws.Range("O17").Value = 0.5797382655268795
$ws.Range("P17").Value = 0.6587494728716036
$ws.Range("Q17").Value = 247.6669650752628
$ws.Range("R17").Value = 1486.001790451577
$ws.Range("S17").Value = 0.2086035231430524
$ws.Range("T17").Value = 0.179558929552401
$ws.Range("G18").Value = 12.3983215
$ws.Range("H18").Value = 24.796643
$ws.Range("I18").Value = 0.3598236230852706
$ws.Range("J18").Value = 0.2725754432404665
$ws.Range("O18").Value = 0.05975306730825326
$ws.Range("P18").Value = 0.06789667671151554
$ws.Range("Q18").Value = 25.52679668423033
$ws.Range("R18").Value = 153.160780105382
$ws.Range("S18").Value = 0.02150056516931373
$ws.Range("T18").Value = 0.01850696674919601
$ws.Range("G19").Value = 12.3983215
$ws.Range("H19").Value = 24.796643
$ws.Range("I19").Value = 0.3598236230852706
$ws.Range("J19").Value = 0.2725754432404665
$ws.Range("M19").Value = 0.01061633333333333
$ws.Range("N19").Value = 0.031849
$ws.Range("O19").Value = 0.0003081068291285177
$ws.Range("P19").Value = 0.0003500980068860779
$ws.Range("Q19").Value = 0.1316247138178333
$ws.Range("R19").Value = 0.789748282907
$ws.Range("S19").Value = 0.0001108641155543376
$ws.Range("T19").Value = 0.00009542811940457657
$ws.Range("G20").Value = 12.3983215
$ws.Range("H20").Value = 24.796643
$ws.Range("I20").Value = 0.3598236230852706
$ws.Range("J20").Value = 0.2725754432404665
$ws.Range("M20").Value = 12.3983215
$ws.Range("N20").Value = 24.796643
$ws.Range("O20").Value = 0.3598236230852706
$ws.Range("P20").Value = 0.2725754432404665
$ws.Range("Q20").Value = 153.7183760173622
$ws.Range("R20").Value = 614.8735040694489
$ws.Range("S20").Value = 0.1294730397302109
$ws.Range("T20").Value = 0.07429737225773675
$ws.Range("G21").Value = 12.3983215
$ws.Range("H21").Value = 24.796643
$ws.Range("I21").Value = 0.3598236230852706
$ws.Range("J21").Value = 0.2725754432404665
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.012988
$ws.Range("N21").Value = 0.038964
$ws.Range("O21").Value = 0.0003769372504682585
$ws.Range("P21").Value = 0.0004283091695283727
$ws.Range("Q21").Value = 0.161029399642
$ws.Range("R21").Value = 0.9661763978519999
$ws.Range("S21").Value = 0.0001356309271392889
$ws.Range("T21").Value = 0.0001167465617281523
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.012988
$ws.Range("H22").Value = 0.038964
$ws.Range("I22").Value = 0.0003769372504682585
$ws.Range("J22").Value = 0.0004283091695283727
$ws.Range("M22").Value = 19.97584633333333
$ws.Range("N22").Value = 59.927539
$ws.Range("O22").Value = 0.5797382655268795
$ws.Range("P22").Value = 0.6587494728716036
$ws.Range("Q22").Value = 0.2594462921773333
$ws.Range("R22").Value = 2.335016629596
$ws.Range("S22").Value = 0.0002185249477989391
$ws.Range("T22").Value = 0.0002821484396528898
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.012988
$ws.Range("H23").Value = 0.038964
$ws.Range("I23").Value = 0.0003769372504682585
$ws.Range("J23").Value = 0.0004283091695283727
$ws.Range("O23").Value = 0.05975306730825326
$ws.Range("P23").Value = 0.06789667671151554
$ws.Range("Q23").Value = 0.02674088063733333
$ws.Range("R23").Value = 0.240667925736
$ws.Range("S23").Value = 0.00002252315689821777
$ws.Range("T23").Value = 0.00002908076921604562
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.012988
$ws.Range("H24").Value = 0.038964
$ws.Range("I24").Value = 0.0003769372504682585
$ws.Range("J24").Value = 0.0004283091695283727
$ws.Range("M24").Value = 0.01061633333333333
$ws.Range("N24").Value = 0.031849
$ws.Range("O24").Value = 0.0003081068291285177
$ws.Range("P24").Value = 0.0003500980068860779
$ws.Range("Q24").Value = 0.0001378849373333333
$ws.Range("R24").Value = 0.001240964436
$ws.Range("S24").Value = 0.000000116136941022197
$ws.Range("T24").Value = 0.0000001499501865829145
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.012988
$ws.Range("H25").Value = 0.038964
$ws.Range("I25").Value = 0.0003769372504682585
$ws.Range("J25").Value = 0.0004283091695283727
$ws.Range("M25").Value = 12.3983215
$ws.Range("N25").Value = 24.796643
$ws.Range("O25").Value = 0.3598236230852706
$ws.Range("P25").Value = 0.2725754432404665
$ws.Range("Q25").Value = 0.161029399642
$ws.Range("R25").Value = 0.9661763978519999
$ws.Range("S25").Value = 0.0001356309271392889
$ws.Range("T25").Value = 0.0001167465617281523
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.012988
$ws.Range("H26").Value = 0.038964
$ws.Range("I26").Value = 0.0003769372504682585
$ws.Range("J26").Value = 0.0004283091695283727
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.012988
$ws.Range("N26").Value = 0.038964
$ws.Range("O26").Value = 0.0003769372504682585
$ws.Range("P26").Value = 0.0004283091695283727
$ws.Range("Q26").Value = 0.000168688144
$ws.Range("R26").Value = 0.001518193296
$ws.Range("S26").Value = 0.0000001420816907905706
$ws.Range("T26").Value = 0.0000001834487447020843

Write-Output "Applied 298 cell updates"
